$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 98
$ws.Cells.Item($row, 1).Value = 46047
$ws.Cells.Item($row, 2).Value = 228
$ws.Cells.Item($row, 3).Value = 230
$ws.Cells.Item($row, 4).Value = 225

$ws.Cells.Item($row, 1).NumberFormat = $ws.Cells.Item($row - 1, 1).NumberFormat
